# feat: added generate strategy logic
#
# 1) Collapse the "{{" / "fundName" / "}}" run-triplet in the fund-name
#    placeholder textbox (slide 1, shape 1) into a single run that carries
#    the (unfilled) formatting of the middle run, i.e. <a:rPr lang="en"
#    sz="1600"/> with text "{{fundName}}".
# 2) Swap the deck's active colour scheme (theme2.xml, the theme used by
#    the slide master / presentation) from "Simple Light" to the
#    "Default" palette.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# --- 1. Merge the {{ / fundName / }} runs -------------------------------
# Runs are currently: [1]="{{" (black solidFill), [2]="fundName" (no fill),
# [3]="}}" (black solidFill). Blank out the outer two so only the
# unfilled middle run survives, then re-expand its text to the full
# placeholder token - this keeps run [2]'s formatting (no solidFill).
$r1 = $tr.Runs(1)
$r1.Text = ""

$tr = $sh.TextFrame.TextRange
$r3 = $tr.Runs(2)
$r3.Text = ""

$tr = $sh.TextFrame.TextRange
$final = $tr.Runs(1)
$final.Text = "{{fundName}}"
$final.Font.Size = 16

# --- 2. Swap the active theme colours from "Simple Light" to "Default" --
$tcs = $s.ThemeColorScheme
$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 5800213    # dk2      158158
$tcs.Item(4).RGB  = 15987699   # lt2      F3F3F3
$tcs.Item(5).RGB  = 13077765   # accent1  058DC7
$tcs.Item(6).RGB  = 3322960    # accent2  50B432
$tcs.Item(7).RGB  = 1791725    # accent3  ED561B
$tcs.Item(8).RGB  = 61421      # accent4  EDEF00
$tcs.Item(9).RGB  = 15059748   # accent5  24CBE5
$tcs.Item(10).RGB = 7529828    # accent6  64E572
$tcs.Item(11).RGB = 13369378   # hlink    2200CC
$tcs.Item(12).RGB = 9116245    # folHlink 551A8B
